# Remove the slide that presented the block-diagram contents as a
# numbered list (the slide immediately before the actual "BLOCK DIAGRAM"
# diagram slide). This is slide index 5 (1-based) in the original deck.
$p = $ppt.ActivePresentation
$p.Slides.Item(5).Delete()
